$wb = $excel.ActiveWorkbook

$oldGuid = "54996d78-d739-42bf-b0be-67f9de45fbc4"
$newGuid = "535a55da-d60b-4d79-be17-5eecfc3feb5b"
$oldHash = "6affe76b585b0e0d0c841b73536a018360eeaa0a"
$newHash = "ff1a53c887d99f4ff0d2bea0e95ac8f7a6612e1a"

# Overview sheet
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("D2").Value = "2016-53-12 22:53:05"

# zh-cn sheet
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-12 22:53:01"

# de-de sheet
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-12 22:53:05"
